$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Hunk 0
$ws.Range("H40").Value = 3231.1667
$ws.Range("I40").Value = 2833.3333
$ws.Range("J40").Value = 3629
$ws.Range("K40").Value = 2833.3333
$ws.Range("L40").Value = 3629
$ws.Range("M40").Value = -2658.3333
$ws.Range("N40").Value = -3979

# Hunk 1
$ws.Range("H43").Value = 37049470
$ws.Range("I43").Value = 83334450
$ws.Range("J43").Value = 21490
$ws.Range("K43").Value = 83334450
$ws.Range("L43").Value = 21490
$ws.Range("M43").Value = -83334381
$ws.Range("N43").Value = -21628

# Hunk 2
$ws.Range("H80").Value = 25000664
$ws.Range("J80").Value = 27778438
$ws.Range("L80").Value = 83335314
$ws.Range("N80").Value = -83337310

# Hunk 3
$ws.Range("H83").Value = 25000664
$ws.Range("J83").Value = 27778438
$ws.Range("L83").Value = 250005942
$ws.Range("N83").Value = -250015926

# Hunk 4
$ws.Range("H98").Value = 3866.7334
$ws.Range("I98").Value = 2999.4443
$ws.Range("J98").Value = 5167.6665
$ws.Range("K98").Value = 2999.4443
$ws.Range("L98").Value = 5167.6665
$ws.Range("M98").Value = -1501.4443
$ws.Range("N98").Value = -8163.6665

# Hunk 5
$ws.Range("H100").Value = 2048.25
$ws.Range("I100").Value = 2331
$ws.Range("K100").Value = 2331
$ws.Range("M100").Value = -1790

# Hunk 6
$ws.Range("H101").Value = 3200
$ws.Range("I101").Value = 2666.6667
$ws.Range("J101").Value = 4000
$ws.Range("K101").Value = 8000.000100000001
$ws.Range("L101").Value = 12000
$ws.Range("M101").Value = -6378.000100000001
$ws.Range("N101").Value = -15244

# Hunk 7
$ws.Range("H107").Value = 23812094
$ws.Range("I107").Value = 2346.4
$ws.Range("K107").Value = 2346.4
$ws.Range("M107").Value = -426.4000000000001

# Hunk 8
$ws.Range("H116").Value = 5997.2856
$ws.Range("I116").Value = 6129.6665
$ws.Range("K116").Value = 6129.6665
$ws.Range("M116").Value = -2687.6665

# Hunk 9
$ws.Range("H122").Value = 3866.7334
$ws.Range("I122").Value = 2999.4443
$ws.Range("J122").Value = 5167.6665
$ws.Range("K122").Value = 8998.332900000001
$ws.Range("L122").Value = 15502.9995
$ws.Range("M122").Value = -6548.332900000001
$ws.Range("N122").Value = -20402.9995

# Hunk 10
$ws.Range("H132").Value = 1629.6757
$ws.Range("I132").Value = 1463.6364
$ws.Range("K132").Value = 4390.9092
$ws.Range("M132").Value = -1860.9092

$ws = $wb.Worksheets.Item("ARM")
# Hunk 11
$ws.Range("H36").Value = 6426.7144
$ws.Range("I36").Value = 6152.077
$ws.Range("J36").Value = 9997
$ws.Range("K36").Value = 6152.077
$ws.Range("L36").Value = 9997
$ws.Range("M36").Value = -5806.077
$ws.Range("N36").Value = -10689

# Hunk 12
$ws.Range("H39").Value = 8335005.5
$ws.Range("I39").Value = 8335005.5
$ws.Range("K39").Value = 8335005.5
$ws.Range("M39").Value = -8334485.5

# Hunk 13
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").ClearContents()
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = 0

# Hunk 14
$ws.Range("H122").Value = 2086.9656
$ws.Range("I122").Value = 2196.8845
$ws.Range("J122").Value = 1134.3334
$ws.Range("K122").Value = 6590.6535
$ws.Range("L122").Value = 3403.0002
$ws.Range("M122").Value = -4140.6535
$ws.Range("N122").Value = -8303.0002

# Hunk 15
$ws.Range("H132").Value = 2247.25
$ws.Range("I132").Value = 973.7727
$ws.Range("K132").Value = 2921.3181
$ws.Range("M132").Value = -391.3181

$ws = $wb.Worksheets.Item("BSM")
# Hunk 16
$ws.Range("H99").Value = 8870.1875
$ws.Range("I99").Value = 9966
$ws.Range("K99").Value = 9966
$ws.Range("M99").Value = -8468

$ws = $wb.Worksheets.Item("CRP")
# Hunk 17
$ws.Range("H20").Value = 129990
$ws.Range("J20").Value = 129990
$ws.Range("L20").Value = 129990
$ws.Range("N20").Value = -130462

# Hunk 18
$ws.Range("H30").Value = 129990
$ws.Range("J30").Value = 129990
$ws.Range("L30").Value = 129990
$ws.Range("N30").Value = -130172

# Hunk 19
$ws.Range("H107").Value = 1992.3334
$ws.Range("I107").Value = 1992.3334
$ws.Range("K107").Value = 1992.3334
$ws.Range("M107").Value = -72.33339999999998

# Hunk 20
$ws.Range("H122").Value = 1768.0625
$ws.Range("I122").Value = 1768.0625
$ws.Range("K122").Value = 5304.1875
$ws.Range("M122").Value = -2854.1875

# Hunk 21
$ws.Range("H128").Value = 129990
$ws.Range("J128").Value = 129990
$ws.Range("L128").Value = 129990
$ws.Range("N128").Value = -139950

# Hunk 22
$ws.Range("H134").Value = 2414.2222
$ws.Range("J134").Value = 3049
$ws.Range("L134").Value = 9147
$ws.Range("N134").Value = -14217

$ws = $wb.Worksheets.Item("CUL")
# Hunk 23
$ws.Range("H5").Value = 1238.8572
$ws.Range("I5").Value = 1024.6364
$ws.Range("J5").Value = 1474.5
$ws.Range("K5").Value = 3073.9092
$ws.Range("L5").Value = 4423.5
$ws.Range("M5").Value = -2961.9092
$ws.Range("N5").Value = -4647.5

# Hunk 24
$ws.Range("H38").Value = 485.84616
$ws.Range("J38").Value = 686.1667
$ws.Range("L38").Value = 2058.5001
$ws.Range("N38").Value = -2752.5001

# Hunk 25
$ws.Range("H68").Value = 2270.889
$ws.Range("I68").Value = 1725
$ws.Range("K68").Value = 5175
$ws.Range("M68").Value = -4364

# Hunk 26
$ws.Range("H71").Value = 2270.889
$ws.Range("I71").Value = 1725
$ws.Range("K71").Value = 15525
$ws.Range("M71").Value = -11469

# Hunk 27
$ws.Range("H122").Value = 11112682
$ws.Range("J122").Value = 2079.5
$ws.Range("L122").Value = 18715.5
$ws.Range("N122").Value = -23615.5

# Hunk 28
$ws.Range("H134").Value = 2219
$ws.Range("I134").Value = 1675.6111
$ws.Range("K134").Value = 5026.8333
$ws.Range("M134").Value = 43.16669999999976

# Hunk 29
$ws.Range("H135").Value = 1238.8572
$ws.Range("I135").Value = 1024.6364
$ws.Range("J135").Value = 1474.5
$ws.Range("K135").Value = 9221.7276
$ws.Range("L135").Value = 13270.5
$ws.Range("M135").Value = -6686.7276
$ws.Range("N135").Value = -18340.5

$ws = $wb.Worksheets.Item("GSM")
# Hunk 30
$ws.Range("H9").Value = 3496.3333
$ws.Range("I9").Value = 244.5
$ws.Range("J9").Value = 10000
$ws.Range("K9").Value = 244.5
$ws.Range("L9").Value = 10000
$ws.Range("M9").Value = -74.5
$ws.Range("N9").Value = -10340

# Hunk 31
$ws.Range("I80").Value = 127710.69
$ws.Range("K80").Value = 127710.69
$ws.Range("M80").Value = -126712.69

# Hunk 32
$ws.Range("I83").Value = 127710.69
$ws.Range("K83").Value = 638553.45
$ws.Range("M83").Value = -633561.45

# Hunk 33
$ws.Range("H107").Value = 77588.46000000001
$ws.Range("I107").Value = 200134
$ws.Range("J107").Value = 997.5
$ws.Range("K107").Value = 200134
$ws.Range("L107").Value = 997.5
$ws.Range("M107").Value = -198214
$ws.Range("N107").Value = -4837.5

# Hunk 34
$ws.Range("H122").Value = 1145.5
$ws.Range("I122").Value = 913.1875
$ws.Range("K122").Value = 2739.5625
$ws.Range("M122").Value = -289.5625

$ws = $wb.Worksheets.Item("LTW")
# Hunk 35
$ws.Range("H7").Value = 8200.526
$ws.Range("I7").Value = 2893.4546
$ws.Range("K7").Value = 2893.4546
$ws.Range("M7").Value = -2781.4546

# Hunk 36
$ws.Range("H40").Value = 6177705.5
$ws.Range("I40").Value = 7941050
$ws.Range("K40").Value = 7941050
$ws.Range("M40").Value = -7940914

# Hunk 37
$ws.Range("H126").Value = 8200.526
$ws.Range("I126").Value = 2893.4546
$ws.Range("K126").Value = 8680.363799999999
$ws.Range("M126").Value = -6210.363799999999

# Hunk 38
$ws.Range("H132").Value = 3127.52
$ws.Range("J132").Value = 2831.1428
$ws.Range("L132").Value = 8493.428400000001
$ws.Range("N132").Value = -13553.4284

$ws = $wb.Worksheets.Item("WVR")
# Hunk 39
$ws.Range("H12").Value = 1504900
$ws.Range("I12").Value = 1504900
$ws.Range("K12").Value = 1504900
$ws.Range("M12").Value = -1504758

# Hunk 40
$ws.Range("H49").Value = 29999
$ws.Range("I49").Value = 29999
$ws.Range("K49").Value = 29999
$ws.Range("M49").Value = -29769

# Hunk 41
$ws.Range("H54").Value = 23332.666
$ws.Range("J54").Value = 49998
$ws.Range("L54").Value = 49998
$ws.Range("N54").Value = -51038

Write-Host "applied all changes"
